$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.281134
$ws.Range("H2").Value = 3.843402
$ws.Range("I2").Value = 0.007312702338676299
$ws.Range("J2").Value = 0.007312702338676299
$ws.Range("M2").Value = 255.0443116666667
$ws.Range("N2").Value = 765.132935
$ws.Range("O2").Value = 0.863617428561108
$ws.Range("P2").Value = 0.8636174285611079
$ws.Range("Q2").Value = 326.7459391827633
$ws.Range("R2").Value = 2940.71345264487
$ws.Range("S2").Value = 0.006315377189560426
$ws.Range("T2").Value = 0.006315377189560425
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.281134
$ws.Range("H3").Value = 3.843402
$ws.Range("I3").Value = 0.007312702338676299
$ws.Range("J3").Value = 0.007312702338676299
$ws.Range("M3").Value = 0.8952453333333334
$ws.Range("O3").Value = 0.003031431940796009
$ws.Range("P3").Value = 0.003031431940796009
$ws.Range("Q3").Value = 1.146929234874667
$ws.Range("R3").Value = 10.322363113872
$ws.Range("S3").Value = [double]"2.2167959442997E-05"
$ws.Range("T3").Value = [double]"2.2167959442997E-05"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.281134
$ws.Range("H4").Value = 3.843402
$ws.Range("I4").Value = 0.007312702338676299
$ws.Range("J4").Value = 0.007312702338676299
$ws.Range("M4").Value = 7.050555333333333
$ws.Range("N4").Value = 21.151666
$ws.Range("O4").Value = 0.02387421396349043
$ws.Range("P4").Value = 0.02387421396349043
$ws.Range("Q4").Value = 9.032706156414665
$ws.Range("R4").Value = 81.294355407732
$ws.Range("S4").Value = 0.0001745850202848748
$ws.Range("T4").Value = 0.0001745850202848748
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.281134
$ws.Range("H5").Value = 3.843402
$ws.Range("I5").Value = 0.007312702338676299
$ws.Range("J5").Value = 0.007312702338676299
$ws.Range("M5").Value = 32.33082866666666
$ws.Range("N5").Value = 96.99248599999999
$ws.Range("O5").Value = 0.1094769255346056
$ws.Range("P5").Value = 0.1094769255346056
$ws.Range("Q5").Value = 41.42012385304132
$ws.Range("R5").Value = 372.781114677372
$ws.Range("S5").Value = 0.0008005721693880017
$ws.Range("T5").Value = 0.0008005721693880016
$ws.Range("I6").Value = 0.9398544320918915
$ws.Range("J6").Value = 0.9398544320918915
$ws.Range("M6").Value = 255.0443116666667
$ws.Range("N6").Value = 765.132935
$ws.Range("O6").Value = 0.863617428561108
$ws.Range("P6").Value = 0.8636174285611079
$ws.Range("Q6").Value = 41994.54659664651
$ws.Range("R6").Value = 377950.9193698186
$ws.Range("S6").Value = 0.8116746678649598
$ws.Range("T6").Value = 0.8116746678649597
$ws.Range("I7").Value = 0.9398544320918915
$ws.Range("J7").Value = 0.9398544320918915
$ws.Range("M7").Value = 0.8952453333333334
$ws.Range("O7").Value = 0.003031431940796009
$ws.Range("P7").Value = 0.003031431940796009
$ws.Range("S7").Value = 0.002849104745142053
$ws.Range("T7").Value = 0.002849104745142053
$ws.Range("I8").Value = 0.9398544320918915
$ws.Range("J8").Value = 0.9398544320918915
$ws.Range("M8").Value = 7.050555333333333
$ws.Range("N8").Value = 21.151666
$ws.Range("O8").Value = 0.02387421396349043
$ws.Range("P8").Value = 0.02387421396349043
$ws.Range("Q8").Value = 1160.915421100915
$ws.Range("R8").Value = 10448.23878990823
$ws.Range("S8").Value = 0.0224382858062966
$ws.Range("T8").Value = 0.0224382858062966
$ws.Range("I9").Value = 0.9398544320918915
$ws.Range("J9").Value = 0.9398544320918915
$ws.Range("M9").Value = 32.33082866666666
$ws.Range("N9").Value = 96.99248599999999
$ws.Range("O9").Value = 0.1094769255346056
$ws.Range("P9").Value = 0.1094769255346056
$ws.Range("Q9").Value = 5323.46117456254
$ws.Range("R9").Value = 47911.15057106286
$ws.Range("S9").Value = 0.1028923736754931
$ws.Range("T9").Value = 0.1028923736754931
$ws.Range("G10").Value = 9.213772333333333
$ws.Range("H10").Value = 27.641317
$ws.Range("I10").Value = 0.05259213672418158
$ws.Range("J10").Value = 0.05259213672418158
$ws.Range("M10").Value = 255.0443116666667
$ws.Range("N10").Value = 765.132935
$ws.Range("O10").Value = 0.863617428561108
$ws.Range("P10").Value = 0.8636174285611079
$ws.Range("Q10").Value = 2349.920222608377
$ws.Range("R10").Value = 21149.28200347539
$ws.Range("S10").Value = 0.04541948588027191
$ws.Range("T10").Value = 0.04541948588027191
$ws.Range("G11").Value = 9.213772333333333
$ws.Range("H11").Value = 27.641317
$ws.Range("I11").Value = 0.05259213672418158
$ws.Range("J11").Value = 0.05259213672418158
$ws.Range("M11").Value = 0.8952453333333334
$ws.Range("O11").Value = 0.003031431940796009
$ws.Range("P11").Value = 0.003031431940796009
$ws.Range("Q11").Value = 8.248586683812444
$ws.Range("R11").Value = 74.23728015431202
$ws.Range("S11").Value = 0.0001594294831003948
$ws.Range("T11").Value = 0.0001594294831003948
$ws.Range("G12").Value = 9.213772333333333
$ws.Range("H12").Value = 27.641317
$ws.Range("I12").Value = 0.05259213672418158
$ws.Range("J12").Value = 0.05259213672418158
$ws.Range("M12").Value = 7.050555333333333
$ws.Range("N12").Value = 21.151666
$ws.Range("O12").Value = 0.02387421396349043
$ws.Range("P12").Value = 0.02387421396349043
$ws.Range("Q12").Value = 64.96221166490244
$ws.Range("R12").Value = 584.659904984122
$ws.Range("S12").Value = 0.001255595924950254
$ws.Range("T12").Value = 0.001255595924950254
$ws.Range("G13").Value = 9.213772333333333
$ws.Range("H13").Value = 27.641317
$ws.Range("I13").Value = 0.05259213672418158
$ws.Range("J13").Value = 0.05259213672418158
$ws.Range("M13").Value = 32.33082866666666
$ws.Range("N13").Value = 96.99248599999999
$ws.Range("O13").Value = 0.1094769255346056
$ws.Range("P13").Value = 0.1094769255346056
$ws.Range("Q13").Value = 297.8888946826735
$ws.Range("R13").Value = 2681.000052144062
$ws.Range("S13").Value = 0.005757625435859026
$ws.Range("T13").Value = 0.005757625435859025
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.042174
$ws.Range("H14").Value = 0.126522
$ws.Range("I14").Value = 0.0002407288452506406
$ws.Range("J14").Value = 0.0002407288452506406
$ws.Range("M14").Value = 255.0443116666667
$ws.Range("N14").Value = 765.132935
$ws.Range("O14").Value = 0.863617428561108
$ws.Range("P14").Value = 0.8636174285611079
$ws.Range("Q14").Value = 10.75623880023
$ws.Range("R14").Value = 96.80614920206999
$ws.Range("S14").Value = 0.0002078976263158431
$ws.Range("T14").Value = 0.0002078976263158431
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.042174
$ws.Range("H15").Value = 0.126522
$ws.Range("I15").Value = 0.0002407288452506406
$ws.Range("J15").Value = 0.0002407288452506406
$ws.Range("M15").Value = 0.8952453333333334
$ws.Range("O15").Value = 0.003031431940796009
$ws.Range("P15").Value = 0.003031431940796009
$ws.Range("Q15").Value = 0.037756076688
$ws.Range("R15").Value = 0.339804690192
$ws.Range("S15").Value = [double]"7.297531105637315E-07"
$ws.Range("T15").Value = [double]"7.297531105637315E-07"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.042174
$ws.Range("H16").Value = 0.126522
$ws.Range("I16").Value = 0.0002407288452506406
$ws.Range("J16").Value = 0.0002407288452506406
$ws.Range("M16").Value = 7.050555333333333
$ws.Range("N16").Value = 21.151666
$ws.Range("O16").Value = 0.02387421396349043
$ws.Range("P16").Value = 0.02387421396349043
$ws.Range("Q16").Value = 0.2973501206279999
$ws.Range("R16").Value = 2.676151085652
$ws.Range("S16").Value = [double]"5.747211958697772E-06"
$ws.Range("T16").Value = [double]"5.747211958697771E-06"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.042174
$ws.Range("H17").Value = 0.126522
$ws.Range("I17").Value = 0.0002407288452506406
$ws.Range("J17").Value = 0.0002407288452506406
$ws.Range("M17").Value = 32.33082866666666
$ws.Range("N17").Value = 96.99248599999999
$ws.Range("O17").Value = 0.1094769255346056
$ws.Range("P17").Value = 0.1094769255346056
$ws.Range("Q17").Value = 1.363520368188
$ws.Range("R17").Value = 12.271683313692
$ws.Range("S17").Value = [double]"2.635425386553599E-05"
$ws.Range("T17").Value = [double]"2.635425386553598E-05"
